# Generate Report for Handoff
# Adds two new "Ready for handoff" rows (60a36207-... and 6461ae72-...)
# to the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$srcRepoCommit = "https://github.com/OpenLocalizationTestOrg/oltest/blob/823e8c08768608709941cae44168b4fd61b662e3/e2e"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Row 4 - 60a36207-1dd1-44c6-a718-743631a49777.md
$ws.Range("A4").Value = "60a36207-1dd1-44c6-a718-743631a49777.md"
$ws.Range("B4").Value = "e2e\60a36207-1dd1-44c6-a718-743631a49777.md"
$ws.Hyperlinks.Add($ws.Range("B4"), "$srcRepoCommit/60a36207-1dd1-44c6-a718-743631a49777.md", "", "", "e2e\60a36207-1dd1-44c6-a718-743631a49777.md")
$ws.Range("C4").Value = ".md"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-08-13 20:32:28"

# Row 5 - 6461ae72-408f-4406-9e3f-a5dbc1c67802.md
$ws.Range("A5").Value = "6461ae72-408f-4406-9e3f-a5dbc1c67802.md"
$ws.Range("B5").Value = "e2e\6461ae72-408f-4406-9e3f-a5dbc1c67802.md"
$ws.Hyperlinks.Add($ws.Range("B5"), "$srcRepoCommit/6461ae72-408f-4406-9e3f-a5dbc1c67802.md", "", "", "e2e\6461ae72-408f-4406-9e3f-a5dbc1c67802.md")
$ws.Range("C5").Value = ".md"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Ready for handoff"
$ws.Range("F5").Value = "Ready for handoff"
$ws.Range("G5").Value = "2016-08-13 20:32:28"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

# Row 4 - 60a36207-1dd1-44c6-a718-743631a49777.md
$ws.Range("A4").Value = "60a36207-1dd1-44c6-a718-743631a49777.md"
$ws.Hyperlinks.Add($ws.Range("A4"), "$srcRepoCommit/60a36207-1dd1-44c6-a718-743631a49777.md", "", "", "60a36207-1dd1-44c6-a718-743631a49777.md")
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "e2e"
$ws.Range("E4").Value = "ht"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "60a36207-1dd1-44c6-a718-743631a49777.3bc2513afad13d82e4d882ae1a8cf4b729dae55d.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-13 20:32:20"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = "0001-01-01 00:00:00"
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = "'False"
$ws.Range("P4").Value = ""

# Row 5 - 6461ae72-408f-4406-9e3f-a5dbc1c67802.md
$ws.Range("A5").Value = "6461ae72-408f-4406-9e3f-a5dbc1c67802.md"
$ws.Hyperlinks.Add($ws.Range("A5"), "$srcRepoCommit/6461ae72-408f-4406-9e3f-a5dbc1c67802.md", "", "", "6461ae72-408f-4406-9e3f-a5dbc1c67802.md")
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "e2e"
$ws.Range("E5").Value = "ht"
$ws.Range("F5").Value = "'False"
$ws.Range("G5").Value = "6461ae72-408f-4406-9e3f-a5dbc1c67802.f0b1cd00a99e22a6a2fdfa6b4e96f6ca53416f9b.zh-cn.xlf"
$ws.Range("H5").Value = "2016-08-13 20:32:20"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "0001-01-01 00:00:00"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "'True"
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "'False"
$ws.Range("P5").Value = ""

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

# Row 4 - 60a36207-1dd1-44c6-a718-743631a49777.md
$ws.Range("A4").Value = "60a36207-1dd1-44c6-a718-743631a49777.md"
$ws.Hyperlinks.Add($ws.Range("A4"), "$srcRepoCommit/60a36207-1dd1-44c6-a718-743631a49777.md", "", "", "60a36207-1dd1-44c6-a718-743631a49777.md")
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "e2e"
$ws.Range("E4").Value = "ht"
$ws.Range("F4").Value = "'False"
$ws.Range("G4").Value = "60a36207-1dd1-44c6-a718-743631a49777.3bc2513afad13d82e4d882ae1a8cf4b729dae55d.de-de.xlf"
$ws.Range("H4").Value = "2016-08-13 20:32:28"
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = "0001-01-01 00:00:00"
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = "'True"
$ws.Range("N4").Value = ""
$ws.Range("O4").Value = "'False"
$ws.Range("P4").Value = ""

# Row 5 - 6461ae72-408f-4406-9e3f-a5dbc1c67802.md
$ws.Range("A5").Value = "6461ae72-408f-4406-9e3f-a5dbc1c67802.md"
$ws.Hyperlinks.Add($ws.Range("A5"), "$srcRepoCommit/6461ae72-408f-4406-9e3f-a5dbc1c67802.md", "", "", "6461ae72-408f-4406-9e3f-a5dbc1c67802.md")
$ws.Range("B5").Value = ".md"
$ws.Range("C5").Value = "Ready for handoff"
$ws.Range("D5").Value = "e2e"
$ws.Range("E5").Value = "ht"
$ws.Range("F5").Value = "'False"
$ws.Range("G5").Value = "6461ae72-408f-4406-9e3f-a5dbc1c67802.f0b1cd00a99e22a6a2fdfa6b4e96f6ca53416f9b.de-de.xlf"
$ws.Range("H5").Value = "2016-08-13 20:32:28"
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = "0001-01-01 00:00:00"
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = "'True"
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "'False"
$ws.Range("P5").Value = ""

# ---------------------------------------------------------------------
# Resize the tables / auto-filters to cover the two newly added rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.ListObjects.Item(1).Resize($wsOverview.Range("A1:G5"))

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.ListObjects.Item(1).Resize($wsZhCn.Range("A1:P5"))

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.ListObjects.Item(1).Resize($wsDeDe.Range("A1:P5"))

Write-Host "Report for handoff generated."
